# Splits the single "Martin Aleksandrov" run into two properly
# spell-checked names ("Valentin Mladenov" and "Aleksanda Angelov"),
# each spread across multiple runs wrapped in <w:proofErr> markers,
# with "Aleksanda Angelov" living in a brand-new paragraph that
# inherits the same list/number formatting.

$d = $word.ActiveDocument

# Locate the paragraph that currently holds "Martin Aleksandrov" (the
# 3rd paragraph in the document: title, "Ivaylo Georgiev", then this one).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Martin Aleksandrov*") {
        $target = $para
        break
    }
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Shared run/paragraph formatting (sz 30, szCs 30, lang en-US) used by
# every run in this list, exactly as the existing runs already have it.
$runProps = '<w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr>'
$paraProps = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
    '<w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr></w:pPr>'

# First paragraph: reuses the original paragraph's own rsid attributes
# so that part of the tag stays untouched; only its content changes.
$p1Attrs = 'w:rsidR="00720316" w:rsidRPr="00103139" w:rsidRDefault="00720316" w:rsidP="00103139"'

$p1 = '<w:p ' + $wNs + ' ' + $p1Attrs + '>' + $paraProps + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $runProps + '<w:t>Valentin</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $runProps + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $runProps + '<w:t>Mladenov</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'

# Second (new) paragraph: brand new, no rsid attributes, carries the
# _GoBack bookmark that used to sit at the end of the original paragraph.
$p2 = '<w:p ' + $wNs + '>' + $paraProps + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $runProps + '<w:t>Aleksanda</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $runProps + '<w:t xml:space="preserve"> Angelov</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

$payload = $p1 + $p2

# Replace the whole paragraph (including its end-of-paragraph mark) with
# the two new paragraphs built above.
$fullRange = $d.Range($target.Range.Start, $target.Range.End)
$fullRange.InsertXML($payload)
